$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @{Row=2; B=25.92709912656175; C=10.46640614322942; E=9.690440854064001; F=60.52218174649737; G=3.776161429866094; J=12.35637135367692; L=10.48919001606729}
    @{Row=3; B=25.76685550090782; C=10.00813087818062; E=9.62434485532677; F=59.63373711618225; G=3.781500622069504; J=12.26654864344047; L=10.57063004378406}
    @{Row=4; B=25.67886903946118; C=9.720251365292899; E=9.584136781358206; F=59.08920117396635; G=3.784942035796313; J=12.21121217109789; L=10.62321181959187}
    @{Row=5; B=25.64566272979932; C=9.60151180564343; E=9.567850365962054; F=58.86772648660101; G=3.786385653456095; J=12.18862039597485; L=10.64529074638652}
    @{Row=6; B=25.64030968669421; C=9.581715585729622; E=9.565152102260225; F=58.83098197964362; G=3.786627859939654; J=12.18486666520723; L=10.64899637805162}
    @{Row=7; B=25.67841044287499; C=9.718655485201101; E=9.583916731877064; F=59.08621231381922; G=3.78496133779874; J=12.21090765325582; L=10.62350694165561}
    @{Row=8; B=25.86970909904969; C=10.30985807611112; E=9.667575794665911; F=60.21575451559102; G=3.777968643118118; J=12.32544180541717; L=10.5167376462807}
    @{Row=9; B=26.32570131257759; C=11.40976305340086; E=9.834406503664155; F=62.43043978886872; G=3.765541484009573; J=12.54836046627907; L=10.3276606639287}
    @{Row=10; B=26.70740593627693; C=12.17264097497655; E=9.958388759157573; F=64.04666737106805; G=3.757182499191449; J=12.71082146175317; L=10.20090740374856}
    @{Row=11; B=26.89054670668018; C=12.50838613377193; E=10.01503396867966; F=64.77734933668825; G=3.753544574254557; J=12.78437817733774; L=10.14583927296055}
    @{Row=12; B=26.96120428076048; C=12.63380015777859; E=10.03651342788057; F=65.0532053928667; G=3.752190449600876; J=12.81217631332395; L=10.12535562765029}
    @{Row=13; B=26.94592980541427; C=12.60686819388495; E=10.03188624424073; F=64.9938349544618; G=3.752481043442464; J=12.80619208648328; L=10.12975076861892}
    @{Row=14; B=26.89633392572275; C=12.51873912695766; E=10.01680054795918; F=64.80006199187116; G=3.753432700168936; J=12.78666630882141; L=10.144146685785}
    @{Row=15; B=26.86612328002477; C=12.46453008410679; E=10.00756373598604; F=64.681256066073; G=3.754018669243731; J=12.774698702125; L=10.1530126159822}
    @{Row=16; B=26.69562337544181; C=12.15046296887772; E=9.954691275494678; F=63.9988082675625; G=3.757423542172871; J=12.70600683381162; L=10.20455811861002}
    @{Row=17; B=26.59341892066999; C=11.95482344940291; E=9.922314214760068; F=63.5788515919747; G=3.759554346835722; J=12.66377317774604; L=10.23684125251391}
    @{Row=18; B=26.53553175050553; C=11.84123893168607; E=9.903715227200143; F=63.33688617856724; G=3.760795435656674; J=12.63944915565702; L=10.25565388198597}
    @{Row=19; B=26.51608818685972; C=11.80260283872738; E=9.897422138778444; F=63.25489482805586; G=3.761218316117383; J=12.63120809307517; L=10.26206556274428}
    @{Row=20; B=26.60420623273745; C=11.97575987897937; E=9.92575842568937; F=63.62360124651464; G=3.759325915595312; J=12.66827241432042; L=10.2333794057076}
    @{Row=21; B=26.91086648331778; C=12.54467232927359; E=10.02123084085888; F=64.85700196609004; G=3.753152539862473; J=12.79240307808248; L=10.13990825558435}
    @{Row=22; B=27.11887003103891; C=12.90639614884257; E=10.08379455593495; F=65.65813912064232; G=3.749254642847602; J=12.87319915763812; L=10.08097139823481}
    @{Row=23; B=27.00718134145128; C=12.71429051058993; E=10.05038989465117; F=65.23107017080054; G=3.751322574263907; J=12.83010908052838; L=10.11223131233251}
    @{Row=24; B=26.59932656979595; C=11.96629796480785; E=9.924201250696722; F=63.60337155436277; G=3.75942913931861; J=12.66623844396026; L=10.23494371952279}
    @{Row=25; B=26.19393834676043; C=11.11958324638411; E=9.789005207669033; F=61.83245774638697; G=3.768767021886843; J=12.48826458664859; L=10.37665996114102}
)

foreach ($row in $data) {
    $r = $row.Row
    $ws.Cells.Item($r, 2).Value = $row.B   # Column B
    $ws.Cells.Item($r, 3).Value = $row.C   # Column C
    $ws.Cells.Item($r, 5).Value = $row.E   # Column E
    $ws.Cells.Item($r, 6).Value = $row.F   # Column F
    $ws.Cells.Item($r, 7).Value = $row.G   # Column G
    $ws.Cells.Item($r, 10).Value = $row.J  # Column J
    $ws.Cells.Item($r, 12).Value = $row.L  # Column L
}
